# "Tabular Loop generation works #422"
#
# The TestItem_Loop example sheet is removed from the workbook and the
# remaining TestItem_Sequence sheet is renamed to ActivitySequence, which
# becomes the single, active sheet of the workbook.

$wb = $excel.ActiveWorkbook

# Drop the "TestItem_Loop" worksheet entirely.
$loopSheet = $wb.Worksheets.Item("TestItem_Loop")
$loopSheet.Delete()

# Rename the remaining sheet and make sure it is the active/selected one.
$seqSheet = $wb.Worksheets.Item("TestItem_Sequence")
$seqSheet.Name = "ActivitySequence"
$seqSheet.Activate()
